# Scheduled market-data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for specific leve rows across multiple crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 6058
$ws.Range("I34").Value = 6058
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6058
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -5855
$ws.Range("N34").ClearContents()

$ws.Range("H36").Value = 6058
$ws.Range("I36").Value = 6058
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 6058
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -5343
$ws.Range("N36").Value = -5343

$ws.Range("H40").Value = 1823.5862
$ws.Range("I40").Value = 1329.1666
$ws.Range("J40").Value = 2172.5881
$ws.Range("K40").Value = 1329.1666
$ws.Range("L40").Value = 2172.5881
$ws.Range("M40").Value = -1154.1666
$ws.Range("N40").Value = -2522.5881

$ws.Range("H64").Value = 3140.2856
$ws.Range("I64").Value = 2989.0908
$ws.Range("J64").Value = 3306.6
$ws.Range("K64").Value = 2989.0908
$ws.Range("L64").Value = 3306.6
$ws.Range("M64").Value = -2741.0908
$ws.Range("N64").Value = -3802.6

$ws.Range("H67").Value = 3140.2856
$ws.Range("I67").Value = 2989.0908
$ws.Range("J67").Value = 3306.6
$ws.Range("K67").Value = 2989.0908
$ws.Range("L67").Value = 3306.6
$ws.Range("M67").Value = -2131.0908
$ws.Range("N67").Value = -5022.6

$ws.Range("H74").Value = 4246.25
$ws.Range("I74").Value = 4020
$ws.Range("J74").Value = 4623.3335
$ws.Range("K74").Value = 4020
$ws.Range("L74").Value = 4623.3335
$ws.Range("M74").Value = -3084
$ws.Range("N74").Value = -6495.3335

$ws.Range("H76").Value = 3150

$ws.Range("H77").Value = 4246.25
$ws.Range("I77").Value = 4020
$ws.Range("J77").Value = 4623.3335
$ws.Range("K77").Value = 20100
$ws.Range("L77").Value = 23116.6675
$ws.Range("M77").Value = -15420
$ws.Range("N77").Value = -32476.6675

$ws.Range("H79").Value = 3150

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1681.8182
$ws.Range("I6").Value = 866.6667
$ws.Range("J6").Value = 1987.5
$ws.Range("K6").Value = 866.6667
$ws.Range("L6").Value = 1987.5
$ws.Range("M6").Value = -693.6667
$ws.Range("N6").Value = -2333.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630

$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184

$ws.Range("H86").Value = 1114.2858
$ws.Range("I86").Value = 1200
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -77
$ws.Range("N86").Value = -3246

$ws.Range("H89").Value = 1114.2858
$ws.Range("I89").Value = 1200
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -384
$ws.Range("N89").Value = -16232

$ws.Range("H105").Value = 2123.6667
$ws.Range("I105").Value = 1980
$ws.Range("J105").Value = 2195.5
$ws.Range("K105").Value = 1980
$ws.Range("L105").Value = 2195.5
$ws.Range("M105").Value = -233
$ws.Range("N105").Value = -5689.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 8000
$ws.Range("J48").Value = 8000
$ws.Range("L48").Value = 8000
$ws.Range("N48").Value = -8952

$ws.Range("H58").Value = 4998.4194
$ws.Range("I58").Value = 2288.6191
$ws.Range("J58").Value = 10689
$ws.Range("K58").Value = 2288.6191
$ws.Range("L58").Value = 10689
$ws.Range("M58").Value = -2085.6191
$ws.Range("N58").Value = -11095

$ws.Range("H62").Value = 14137.223
$ws.Range("I62").Value = 3033.5715
$ws.Range("K62").Value = 3033.5715
$ws.Range("M62").Value = -2409.5715

$ws.Range("H65").Value = 14137.223
$ws.Range("I65").Value = 3033.5715
$ws.Range("K65").Value = 15167.8575
$ws.Range("M65").Value = -12047.8575

$ws.Range("H136").Value = 4998.4194
$ws.Range("I136").Value = 2288.6191
$ws.Range("J136").Value = 10689
$ws.Range("K136").Value = 6865.8573
$ws.Range("L136").Value = 32067
$ws.Range("M136").Value = -4315.8573
$ws.Range("N136").Value = -37167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 2633.3333
$ws.Range("J42").Value = 2633.3333
$ws.Range("L42").Value = 7899.999899999999
$ws.Range("N42").Value = -8967.999899999999

$ws.Range("H64").Value = 8003.5
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 6000
$ws.Range("M64").Value = -5730

$ws.Range("H67").Value = 8003.5
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 6000
$ws.Range("M67").Value = -5064

$ws.Range("H113").Value = 656.8946999999999
$ws.Range("I113").Value = 523.1429000000001
$ws.Range("J113").Value = 734.9167
$ws.Range("K113").Value = 1569.4287
$ws.Range("L113").Value = 2204.7501
$ws.Range("M113").Value = 600.5712999999998
$ws.Range("N113").Value = -6544.7501

$ws.Range("H132").Value = 1334.25
$ws.Range("I132").Value = 839
$ws.Range("J132").Value = 1559.3636
$ws.Range("K132").Value = 7551
$ws.Range("L132").Value = 14034.2724
$ws.Range("M132").Value = -5021
$ws.Range("N132").Value = -19094.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 470.9
$ws.Range("I13").Value = 463.75
$ws.Range("J13").Value = 499.5
$ws.Range("K13").Value = 463.75
$ws.Range("L13").Value = 499.5
$ws.Range("M13").Value = -324.75
$ws.Range("N13").Value = -777.5

$ws.Range("H70").Value = 6869.3
$ws.Range("I70").Value = 7198.2354
$ws.Range("J70").Value = 6439.154
$ws.Range("K70").Value = 7198.2354
$ws.Range("L70").Value = 6439.154
$ws.Range("M70").Value = -6928.2354
$ws.Range("N70").Value = -6979.154

$ws.Range("H73").Value = 6869.3
$ws.Range("I73").Value = 7198.2354
$ws.Range("J73").Value = 6439.154
$ws.Range("K73").Value = 7198.2354
$ws.Range("L73").Value = 6439.154
$ws.Range("M73").Value = -6262.2354
$ws.Range("N73").Value = -8311.154

$ws.Range("H136").Value = 32442
$ws.Range("J136").Value = 32442
$ws.Range("L136").Value = 97326
$ws.Range("N136").Value = -102426
